$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 610 entirely; all rows below shift up by one.
$ws.Rows.Item(610).Delete()
